$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A81").Value = "GRT-USD"
